$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped from 45178 (2023-09-09)
# to 45179 (2023-09-10) for every data row (rows 2 through 428).
$ws.Range("C2:C428").Value = 45179
